$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("B3").Value = 972
$ws.Range("C3").Value = 1962
$ws.Range("D3").Value = 3711
$ws.Range("E3").Value = 8537
$ws.Range("F3").Value = 12800
$ws.Range("G3").Value = 12400

# Row 8
$ws.Range("B8").Value = 6972
$ws.Range("C8").Value = 13800
$ws.Range("D8").Value = 28100
$ws.Range("E8").Value = 62900
$ws.Range("F8").Value = 102000
$ws.Range("G8").Value = 158000

# Row 13
$ws.Range("B13").Value = 7613
$ws.Range("C13").Value = 10800
$ws.Range("D13").Value = 12800
$ws.Range("E13").Value = 13500
$ws.Range("F13").Value = 12200
$ws.Range("G13").Value = 12500

# Row 18
$ws.Range("B18").Value = 129000
$ws.Range("C18").Value = 179000
$ws.Range("D18").Value = 199000
$ws.Range("E18").Value = 203000
$ws.Range("F18").Value = 178000
$ws.Range("G18").Value = 149000

# Row 23
$ws.Range("B23").Value = 3820
$ws.Range("C23").Value = 3305
$ws.Range("D23").Value = 6787
$ws.Range("E23").Value = 7454
$ws.Range("F23").Value = 10400
$ws.Range("G23").Value = 7043

# Row 28
$ws.Range("B28").Value = 73600
$ws.Range("D28").Value = 239000
$ws.Range("E28").Value = 262000
$ws.Range("F28").Value = 348000
$ws.Range("G28").Value = 236000

# Row 33
$ws.Range("B33").Value = 6224
$ws.Range("C33").Value = 8551
$ws.Range("D33").Value = 8943
$ws.Range("E33").Value = 9959
$ws.Range("F33").Value = 9617
$ws.Range("G33").Value = 11000

# Row 38
$ws.Range("B38").Value = 131000
$ws.Range("C38").Value = 167000
$ws.Range("D38").Value = 178000
$ws.Range("E38").Value = 193000
$ws.Range("F38").Value = 189000
$ws.Range("G38").Value = 176000
